$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Insert a new column at I (this shifts "Lane" and everything after it one
# column to the right, turning the old H column "Completed Tests" into the
# new "Number of Stations" slot and opening up I for the new "Number of
# Drops" column).
# ---------------------------------------------------------------------------
$ws.Columns("I").Insert() | Out-Null

# Rename the (now relocated) "Completed Tests" header to "Number of Stations"
# and give it the same number format as its F/G neighbours (time format).
$ws.Range("H1").Value = "Number of Stations"
$ws.Range("H1").NumberFormat = "h:mm"

# Populate the brand-new column with its header (inherits the plain header
# style already, same as D1/E1/J1/etc.).
$ws.Range("I1").Value = "Number of Drops"

# ---------------------------------------------------------------------------
# Fix up column widths that shifted because of the inserted column.
# ---------------------------------------------------------------------------
# H used to hold the "Completed Tests" header (custom width); now that the
# custom width belongs to the new I column, put H back to the plain
# (F/G-style) width and give I the width that used to belong to H.
$ws.Columns("H").ColumnWidth = 8.307291666666666
$ws.Columns("I").ColumnWidth = 10.877604166666666

# J (old "Insufficient Field Tests" width) -> K
$ws.Columns("K").ColumnWidth = 9.307291666666666

# S:T (old "Decreasing Deflections"/"Station < Section Length" width) -> T:U
$ws.Range("T1:U1").ColumnWidth = 9.592447916666666

# U (old "Comments" width) -> V
$ws.Columns("V").ColumnWidth = 10.451822916666666

# ---------------------------------------------------------------------------
# Conditional formatting ranges all need to slide one column to the right,
# matching the column insertion at I.
# ---------------------------------------------------------------------------
$cfs = $ws.Cells.FormatConditions

for ($i = 1; $i -le $cfs.Count; $i++) {
    $cf = $cfs.Item($i)
    $addr = $cf.AppliesTo.Address()

    if ($addr -eq '$A$1:$BC$1048576') {
        $cf.ModifyAppliesToRange($ws.Range("A1:BD1048576"))
    }
    elseif ($addr -eq '$J$2:$J$1040000') {
        $cf.ModifyAppliesToRange($ws.Range("K2:K1040000,T2:U1040000"))
        $cf.Formula1 = "=LEN(TRIM(K2))>0"
    }
    elseif ($addr -eq '$A$1:$BZ$1048576') {
        $cf.ModifyAppliesToRange($ws.Range("A1:CA1048576"))
    }
}

# ---------------------------------------------------------------------------
# Match the author's final cursor position from the saved file.
# ---------------------------------------------------------------------------
$ws.Range("H2").Select() | Out-Null
